# Updated parameter library and browse product definitions.
# Rewrites the BrowseProduct/Param1/Param2/Param3 table on Sheet1:
#   - HY2  -> renamed to HYD2 (same params)
#   - HY3  -> renamed to HYD3 (same params)
#   - MAF2 row removed
#   - New rows added: FM3, PLG, SUL
#   - Table re-sorted alphabetically by BrowseProduct (column A)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Final data, already sorted alphabetically by BrowseProduct.
$data = @(
  @("BIO", "BD1200",      "BD670",     "D700"),
  @("CR2", "MIN2295_2480","MIN2345_2537","BDCARB"),
  @("FAL", "R2529",       "R1506",     "R1080"),
  @("FM2", "BD530_2",     "BD920_2",   "BDI1000VIS"),
  @("FM3", "BD670",       "BD875",     "BD905"),
  @("HEM", "BD530_2",     "BD670",     "BD875"),
  @("HYD", "SINDEX2",     "BD2100_2",  "BD1900_2"),
  @("HYD2","BD1200",      "BD1450",    "BD1900r2"),
  @("HYD3","BA1200",      "BA1450",    "BA1900"),
  @("HYS", "MIN2250",     "BD2250",    "BD1900r2"),
  @("MAF", "OLINDEX3",    "LCPINDEX2", "HCPINDEX2"),
  @("PAL", "BD2210_2",    "BD2190",    "BD2165"),
  @("PFM", "BD2355",      "D2300",     "BD2290"),
  @("PHY", "D2200",       "D2300",     "BD1900r2"),
  @("PLG", "BD1300",      "RPEAK1",    "LCPINDEX2"),
  @("SED", "BDCARB",      "BD2100_3",  "GINDEX"),
  @("SUL", "GINDEX",      "SINDEX2",   "BD2265"),
  @("TRU", "R637",        "R550",      "R463")
)

$rowCount = $data.Length
$lastRow = 1 + $rowCount   # header is row 1

for ($i = 0; $i -lt $rowCount; $i++) {
  $r = $i + 2
  $rowVals = $data[$i]
  $ws.Cells.Item($r, 1).Value = $rowVals[0]
  $ws.Cells.Item($r, 2).Value = $rowVals[1]
  $ws.Cells.Item($r, 3).Value = $rowVals[2]
  $ws.Cells.Item($r, 4).Value = $rowVals[3]
}

# The old table had 17 rows (through row 17); the new one has 19. No extra
# rows need clearing since we only grew the range, but make sure nothing
# stale is left over below the new table from a previous longer table.
# (Not applicable here since old=17 < new=19.)

# Re-apply the sort (data is already sorted, this captures sort state/metadata,
# matching the workbook's recorded <sortState> for A2:D19 sorted by A2:A19).
$dataRange = $ws.Range("A1:D" + $lastRow)
$keyRange = $ws.Range("A2:A" + $lastRow)
$srt = $ws.Sort
$srt.SortFields.Clear()
$srt.SortFields.Add($keyRange)
$srt.SetRange($dataRange)
$srt.Header = 1
$srt.Apply()

# Update the hidden _FilterDatabase defined name to cover the new range.
foreach ($n in $wb.Names) {
  if ($n.Name -eq "Sheet1!_FilterDatabase") {
    $n.RefersTo = "=Sheet1!`$A`$1:`$D`$" + $lastRow
  }
}

# Move the active selection to B19 (last edited row), matching the saved view.
[void]$ws.Range("B" + $lastRow).Select()

# Nudge the window position, matching the saved workbook view (best effort).
$excel.ActiveWindow.Left = 3360
